$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.079.09"
$ws.Range("E2").Value = "  -8.94%  "

# Row 3
$ws.Range("D3").Value = "2.394.42"
$ws.Range("E3").Value = "  -15.84%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'460.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.07%  "

# Row 6
$ws.Range("D6").Value = "130.41"
$ws.Range("E6").Value = "  -4.79%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").Value = "0.484"
$ws.Range("E8").Value = "  -8.33%  "

# Row 9
$ws.Range("D9").Value = "2.413.93"
$ws.Range("E9").Value = "  -15.14%  "

# Row 10
$ws.Range("E10").Value = "  -9.29%  "

# Row 11
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -10.61%  "

# Row 12
$ws.Range("D12").Value = "0.316"
$ws.Range("E12").Value = "  -9.24%  "

# Row 13
$ws.Range("D13").Value = "0.122"
$ws.Range("E13").Value = "  -4.20%  "

# Row 14
$ws.Range("D14").Value = "2.829.97"
$ws.Range("E14").Value = "  -15.43%  "

# Row 15
$ws.Range("D15").Value = "53.937.51"
$ws.Range("E15").Value = "  -9.25%  "

# Row 16
$ws.Range("D16").Value = "19.56"
$ws.Range("E16").Value = "  -10.34%  "

# Row 17
$ws.Range("E17").Value = "  -5.56%  "

# Row 18
$ws.Range("D18").Value = "2.421.13"
$ws.Range("E18").Value = "  -14.85%  "

# Row 19
$ws.Range("E19").Value = "  -11.94%  "

# Row 20
$ws.Range("D20").Value = "309.88"
$ws.Range("E20").Value = "  -12.25%  "

# Row 21
$ws.Range("D21").Value = "9.36"
$ws.Range("E21").Value = "  -15.91%  "

# Row 22
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("D23").Value = "5.65"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  -15.07%  "

# Row 25
$ws.Range("D25").Value = "'56.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.88%  "

# Row 26
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.155"
$ws.Range("E27").Value = "  -10.74%  "

# Row 28
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").Value = "0.379"
$ws.Range("E28").Value = "  -11.87%  "

# Row 29
$ws.Range("D29").Value = "2.489.46"
$ws.Range("E29").Value = "  -16.46%  "

# Row 30
$ws.Range("D30").Value = "7.11"
$ws.Range("E30").Value = "  -5.40%  "

# Row 31
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.16%  "

# Row 32
$ws.Range("E32").Value = "  -14.18%  "

# Row 33
$ws.Range("D33").Value = "149.53"
$ws.Range("E33").Value = "  -1.01%  "

# Row 34
$ws.Range("D34").Value = "17.61"
$ws.Range("E34").Value = "  -7.78%  "

# Row 35
$ws.Range("D35").Value = "'1.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -13.66%  "

# Row 36
$ws.Range("D36").Value = "5.01"
$ws.Range("E36").Value = "  -6.98%  "

# Row 37
$ws.Range("D37").Value = "3.46"
$ws.Range("E37").Value = "  -17.55%  "

# Row 38
$ws.Range("E38").Value = "  -9.52%  "

# Row 39
$ws.Range("D39").Value = "0.796"
$ws.Range("E39").Value = "  -15.04%  "

# Row 40
$ws.Range("D40").Value = "33.68"
$ws.Range("E40").Value = "  -7.85%  "

# Row 41
$ws.Range("D41").Value = "0.994"
$ws.Range("E41").Value = "  -0.50%  "

# Row 42
$ws.Range("E42").Value = "  -4.54%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0527"
$ws.Range("E43").Value = "  -6.31%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  -7.34%  "

# Row 45
$ws.Range("D45").Value = "10.14"
$ws.Range("E45").Value = "  -1.99%  "

# Row 46
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  -10.89%  "

# Row 47
$ws.Range("D47").Value = "1.960.17"
$ws.Range("E47").Value = "  -12.21%  "

# Row 48
$ws.Range("E48").Value = "  -5.13%  "

# Row 49
$ws.Range("D49").Value = "0.0863"
$ws.Range("E49").Value = "  -3.06%  "

# Row 50
$ws.Range("D50").Value = "'4.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.65%  "

# Row 51
$ws.Range("D51").Value = "16.41"
$ws.Range("E51").Value = "  -16.53%  "
